# Weekly update: insert a new most-recent price record for
# "Terminal La Palmera de La Serena - Poroto granado" and push the
# existing history down by one row (row 43 -> 44, 44 -> 45, ... 65 -> 66).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row right above the current first data-history row
# (row 43). Excel's InsertShiftDown semantics take care of moving every
# row below it (old 43..65) down to (44..66), values, formulas and
# formatting intact.
$ws.Rows.Item(43).Insert()

# Populate the freshly inserted row 43 with this week's record.
$ws.Cells.Item(43, 1).Value  = 8
$ws.Cells.Item(43, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(43, 3).Value  = "Coquimbo"
$ws.Cells.Item(43, 4).Value  = 44596
$ws.Cells.Item(43, 5).Value  = 4
$ws.Cells.Item(43, 6).Value  = 100112030
$ws.Cells.Item(43, 7).Value  = "Poroto granado"
$ws.Cells.Item(43, 8).Value  = "Sin especificar"
$ws.Cells.Item(43, 9).Value  = "Primera"
$ws.Cells.Item(43, 10).Value = 500
$ws.Cells.Item(43, 11).Value = 31000
$ws.Cells.Item(43, 12).Value = 32000
$ws.Cells.Item(43, 13).Value = 31500
$ws.Cells.Item(43, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(43, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(43, 16).Value = 1260
$ws.Cells.Item(43, 17).Value = 25
$ws.Cells.Item(43, 18).Value = "Hortaliza"
